# Update "合肥-漫展信息.xlsx": refresh scraped bilibili show listings.
#
# 1) The expired event "合肥·ACGN夏日游园会预热场" (id=89914) is removed from
#    the "展览" (sheet 1) and "全部类型" (sheet 4) listings, and the
#    remaining rows / index numbers in column A shift up to stay sequential.
# 2) Several events have updated "想去人数" (want-to-go count, column F)
#    values reflecting newer scrape data.

$wb = $excel.ActiveWorkbook

# want-to-go count (column F) updates, keyed by the unique bilibili show id
# found in column H, applied to every sheet that contains that id.
$countUpdates = @{
    "id=84767" = 7614
    "id=90111" = 27
    "id=88650" = 4279
    "id=89420" = 328
    "id=89421" = 596
    "id=88430" = 677
    "id=89803" = 159
    "id=90032" = 10
}

$removeId = "id=89914"

for ($i = 1; $i -le $wb.Worksheets.Count(); $i++) {
    $ws = $wb.Worksheets.Item($i)

    $lastRow = $ws.Cells.Item($ws.Rows.Count(), 1).End(-4162).Row()
    if ($lastRow -lt 2) {
        continue
    }

    # Remove the row for the expired event, if this sheet has it.
    $hRange = $ws.Range($ws.Cells.Item(2, 8), $ws.Cells.Item($lastRow, 8))
    $removeCell = $hRange.Find($removeId)
    if ($removeCell -ne $null) {
        $ws.Rows.Item($removeCell.Row()).Delete()
        $lastRow = $lastRow - 1
    }

    # Renumber the sequential index in column A (row 1 is the header "0",
    # data rows should read 1, 2, 3, ... in order).
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Apply refreshed "想去人数" counts to whichever rows are present.
    if ($lastRow -ge 2) {
        $hRange2 = $ws.Range($ws.Cells.Item(2, 8), $ws.Cells.Item($lastRow, 8))
        foreach ($key in $countUpdates.Keys) {
            $cell = $hRange2.Find($key)
            if ($cell -ne $null) {
                $ws.Cells.Item($cell.Row(), 6).Value = $countUpdates[$key]
            }
        }
    }
}
